# Fix template surat undangan & berkas administrasi
#
# Two table cells in the document each contain a small paragraph whose
# text is literally "1. " (a leftover numbering placeholder). For both
# cells we need to:
#   1. strip the leading "1." off that paragraph so only the trailing
#      space remains;
#   2. add an extra blank paragraph after the (already existing) blank
#      paragraph that follows it;
#   3. add a further, centered paragraph containing the
#      ${validator_administrasi} merge placeholder.
#
# The second table cell gets one more blank paragraph than the first
# (it already had one fewer to start with), matching the target diff.

$d = $word.ActiveDocument

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$paraEmptyXml = "<w:p $ns>" +
    "<w:pPr>" +
        "<w:rPr>" +
            "<w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/>" +
            "<w:noProof/>" +
            "<w:sz w:val=`"16`"/>" +
            "<w:szCs w:val=`"16`"/>" +
        "</w:rPr>" +
    "</w:pPr>" +
"</w:p>"

$paraValidatorXml = "<w:p $ns>" +
    "<w:pPr>" +
        "<w:jc w:val=`"center`"/>" +
        "<w:rPr>" +
            "<w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/>" +
            "<w:noProof/>" +
            "<w:sz w:val=`"16`"/>" +
            "<w:szCs w:val=`"16`"/>" +
        "</w:rPr>" +
    "</w:pPr>" +
    "<w:r>" +
        "<w:rPr>" +
            "<w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/>" +
            "<w:noProof/>" +
            "<w:sz w:val=`"16`"/>" +
            "<w:szCs w:val=`"16`"/>" +
        "</w:rPr>" +
        "<w:t>`${validator_administrasi}</w:t>" +
    "</w:r>" +
"</w:p>"

function Find-NumberingParagraph($doc) {
    # Locate the next remaining paragraph whose whole text is "1." or "1. "
    # followed by the paragraph mark. Re-querying fresh (instead of caching
    # Paragraph objects across InsertXML calls) keeps this reliable even
    # after earlier edits shift the document.
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        if ($t -eq "1.`r" -or $t -eq "1. `r") {
            return $p
        }
    }
    return $null
}

for ($i = 0; $i -lt 2; $i++) {
    $target = Find-NumberingParagraph $d

    # Strip the leading "1." (two characters), leaving only the trailing
    # space (if any) in place with its original run/formatting.
    $head = $d.Range($target.Range.Start, $target.Range.Start + 2)
    $head.Text = ""

    # The paragraph right after the numbering paragraph is an existing
    # blank paragraph — leave it untouched, but append a new blank
    # paragraph after it.
    $afterBlank = $target.Next()
    $insertion1 = $d.Range($afterBlank.Range.End, $afterBlank.Range.End)
    $insertion1.InsertXML($paraEmptyXml)

    # The second location needs one additional blank paragraph compared
    # to the first.
    if ($i -eq 1) {
        $secondBlank = $target.Next().Next()
        $insertion1b = $d.Range($secondBlank.Range.End, $secondBlank.Range.End)
        $insertion1b.InsertXML($paraEmptyXml)
    }

    # Finally, append the centered paragraph carrying the
    # ${validator_administrasi} placeholder.
    $lastBlank = $target.Next()
    for ($k = 0; $k -lt $i; $k++) {
        $lastBlank = $lastBlank.Next()
    }
    $lastBlank = $lastBlank.Next()
    $insertion2 = $d.Range($lastBlank.Range.End, $lastBlank.Range.End)
    $insertion2.InsertXML($paraValidatorXml)
}

Write-Output "Done applying validator_administrasi edits"
